# Apply updated cryptocurrency price/volume data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '29.625.88'
$ws.Cells.Item(2, 5).Value = '  +0.95%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.853.16'
$ws.Cells.Item(3, 5).Value = '  +0.33%  '

# Row 4
$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9989'
$cell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.02%  '

# Row 5
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '241.05'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.27%  '

# Row 6
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.6309'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.32%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  -0.01%  '

# Row 8
$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07493'
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -1.12%  '

# Row 9
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.2917'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  +0.05%  '

# Row 10
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '24.93'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +1.55%  '

# Row 11
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07753'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -0.03%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '1.854.97'
$ws.Cells.Item(12, 5).Value = '  +0.44%  '

# Row 13
$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.044'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +0.59%  '

# Row 14
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.6832'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.72%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  -1.09%  '

# Row 16
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '82.86'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -0.25%  '

# Row 17
$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.305'
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +3.30%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '29.626.26'
$ws.Cells.Item(18, 5).Value = '  +0.87%  '

# Row 19
$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '230.95'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.78%  '

# Row 20
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.41'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.65%  '

# Row 21
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.582'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +1.88%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +0.02%  '

# Row 23
$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.001'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.02%  '

# Row 24
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '159.64'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.32%  '

# Row 25
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.530'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +0.94%  '

# Row 26
$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1370'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -1.86%  '

# Row 27
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '17.60'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -0.41%  '

# Row 28
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.06636'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +16.83%  '

# Row 29
$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.453'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +1.96%  '

# Row 30
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.489'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +1.23%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.122'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  +2.16%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'Filecoin'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.113'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  +0.07%  '

# Row 33
$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.846'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  +1.31%  '

# Row 34
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.146'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -0.78%  '

# Row 35
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.6991'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  +0.47%  '

# Row 36
$cell = $ws.Cells.Item(36, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.567'
$cell.Style = "Normal"
$ws.Cells.Item(36, 5).Value = '  -0.74%  '

# Row 37
$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.01875'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  +2.29%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '1.265.19'
$ws.Cells.Item(38, 5).Value = '  +1.98%  '

# Row 39
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.848'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +4.71%  '

# Row 40
$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.783'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +6.33%  '

# Row 41
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.9337'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +3.54%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(42, 4).Value = '2.032.95'
$ws.Cells.Item(42, 5).Value = '  +1.34%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'PaxDollar'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.002'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +0.23%  '

# Row 44
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '101.48'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +0.05%  '

# Row 45
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '66.39'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +1.35%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +2.81%  '

# Row 47
$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.743'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  +3.93%  '

# Row 48
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.120'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +0.15%  '

# Row 49
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1166'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  +1.40%  '

# Row 50
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '9.016'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +0.20%  '

# Row 51
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.3963'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -0.81%  '

